$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells in column D hold price figures formatted as plain text (e.g. "3.504.19")
# which Excel would otherwise auto-parse as a number. Force text format, assign,
# then restore the default "Normal" style so no stray number format lingers on
# the cell (matching the original workbook, which has no custom formatting there).

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "66.955.86"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.78%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.504.88"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -2.35%  "

$ws.Range("E4").Value = "  +0.04%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "201.82"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +4.97%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "550.61"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -4.52%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.495.21"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -2.48%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.602"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -2.52%  "

$ws.Range("E9").Value = "  +0.20%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.655"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -3.07%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "60.68"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +10.97%  "

$ws.Range("E12").Value = "  -4.73%  "

$ws.Range("E13").Value = "  -0.38%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "9.79"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.22%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "4.076.28"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -2.05%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.503.03"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -2.38%  "

$ws.Range("E17").Value = "  -0.76%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "18.44"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.49%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "66.737.73"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.01%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "11.78"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -4.02%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.03"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -3.53%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "387.92"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -3.45%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "3.99"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -5.01%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.92"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -10.03%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "82.33"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -3.70%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "6.14"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.62%  "

$ws.Range("E27").Value = "  -4.93%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "11.93"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -4.83%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "3.70"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.56%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "8.83"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -3.16%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "30.56"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -2.19%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "7.31"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -9.26%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "672.58"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +1.31%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "11.69"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -4.16%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "63.18"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -1.29%  "

$ws.Range("E36").Value = "  -5.16%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "39.23"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -7.77%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.406"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -2.25%  "

$ws.Range("E39").Value = "  -0.22%  "

$ws.Range("B40").Value = "Maker"
$ws.Range("C40").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.154.65"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.88%  "

$ws.Range("B41").Value = "ThetaToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.07"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -2.56%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.999"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.03%  "

$ws.Range("E43").Value = "  -3.49%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0₃0702"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -10.84%  "

$ws.Range("B45").Value = "Fetch.AI"
$ws.Range("C45").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.54"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -12.18%  "

$ws.Range("B46").Value = "dogwifhat"
$ws.Range("C46").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.79"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +14.93%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.72"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +7.56%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0397"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -4.97%  "

$ws.Range("E49").Value = "  -3.46%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "137.30"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -3.99%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "8.29"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -5.55%  "
